# Add a new "Billing_Address" worksheet after Sheet1, populate it with the
# billing-address table, reuse Sheet1's existing header/body cell styles
# (via copy/paste-special so no new style entries get created), and make
# the new sheet the active tab/selection.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New sheet goes right after Sheet1.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Billing_Address"

# Header row.
$ws2.Range("A1").Value = "City"
$ws2.Range("B1").Value = "Address1"
$ws2.Range("C1").Value = "Zip"
$ws2.Range("D1").Value = "Phone number"

# Data rows.
$ws2.Range("A2").Value = "Chittoor"
$ws2.Range("B2").Value = "Ragigunta"
$ws2.Range("C2").Value = 517600
$ws2.Range("D2").Value = 8908908907

$ws2.Range("A3").Value = "Lanka"
$ws2.Range("B3").Value = "Puram"
$ws2.Range("C3").Value = 600002
$ws2.Range("D3").Value = 9019019012

# Reuse Sheet1's existing header style (yellow fill + border + centered)
# for the new header row, and Sheet1's existing bordered/centered body
# style for the new data rows - via copy/paste-special so the workbook's
# style table isn't expanded with duplicate entries.
$null = $ws1.Range("A1:B1").Copy()
$null = $ws2.Range("A1:D1").PasteSpecial(-4122)

$null = $ws1.Range("B2").Copy()
$null = $ws2.Range("A2:D3").PasteSpecial(-4122)

# Make the new sheet the active tab with its selection on E11.
$null = $ws2.Activate()
$null = $ws2.Range("E11").Select()
